$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = "0130_215757.png"
$ws.Range("I2").Value = "01-30 21:57:57 setText execution was Passed"

$ws.Range("G3").Value = "0130_215758.png"
$ws.Range("I3").Value = "01-30 21:57:58 click execution was Passed"
